$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "C"="0.020209235148819"; "D"="0.2015911232664536"; "E"="0.1694573906331911"; "F"="1.322367260782315"; "G"="0.002444975434858508"; "I"="0.6698032520515866"; "J"="0.1854673713836164"; "K"="1.288168344233554"; "N"="1.15462960317079"; "O"="3.100628091771824" }
    3 = @{ "C"="0.01820130350662907"; "D"="0.196446172594591"; "E"="0.1654361252300305"; "F"="1.320193292982076"; "G"="0.002447927713713043"; "I"="0.6686332824643628"; "J"="0.181294522301009"; "K"="1.155100810727276"; "N"="1.155773944436731"; "O"="3.106970680326697" }
    4 = @{ "C"="0.0169600241453054"; "D"="0.1933618472487524"; "E"="0.1630505293028435"; "F"="1.319683274385774"; "G"="0.002449837522238838"; "I"="0.6683744704613517"; "J"="0.1788407846000126"; "K"="1.073395433905404"; "N"="1.156918539267267"; "O"="3.112931772750329" }
    5 = @{ "C"="0.01645209885411703"; "D"="0.1921238400796028"; "E"="0.1620994056710821"; "F"="1.319682687183175"; "G"="0.00245064027405693"; "I"="0.6683844164218655"; "J"="0.1778681175703269"; "K"="1.04010154388439"; "N"="1.157496389660608"; "O"="3.115879979731488" }
    6 = @{ "C"="0.01636763233231875"; "D"="0.1919194131471471"; "E"="0.1619427432787965"; "F"="1.319695102543385"; "G"="0.002450775051716688"; "I"="0.6683930344378126"; "J"="0.1777082525877631"; "K"="1.034573283143487"; "N"="1.15759907915043"; "O"="3.116400854550619" }
    7 = @{ "C"="0.01695318253896261"; "D"="0.1933450744745642"; "E"="0.1630376169360552"; "F"="1.319682427528562"; "G"="0.002449848249183003"; "I"="0.6683741374824592"; "J"="0.1788275565455493"; "K"="1.07294641118537"; "N"="1.156925880851205"; "O"="3.112969432743341" }
    8 = @{ "C"="0.0195186546218622"; "D"="0.1998016979041211"; "E"="0.1680535519556798"; "F"="1.321446390206319"; "G"="0.002445973272325199"; "I"="0.6693043765805058"; "J"="0.184006065489676"; "K"="1.242288136187483"; "N"="1.154932607818438"; "O"="3.102385652361647" }
    9 = @{ "C"="0.02448228283116549"; "D"="0.2130523031401736"; "E"="0.1785513022097334"; "F"="1.33145917150901"; "G"="0.002439141515904169"; "I"="0.6747826166623696"; "J"="0.1950224772009364"; "K"="1.574281778734075"; "N"="1.154519410901869"; "O"="3.098063831071158" }
    10 = @{ "C"="0.02808759383898973"; "D"="0.223143213431257"; "E"="0.1866671897041741"; "F"="1.342826747285955"; "G"="0.002434585130519534"; "I"="0.6810474044763311"; "J"="0.2036441398949762"; "K"="1.818070370856731"; "N"="1.156333650666426"; "O"="3.104959759751864" }
    11 = @{ "C"="0.02971866015882085"; "D"="0.2278104004347767"; "E"="0.1904469144012495"; "F"="1.348872778335178"; "G"="0.002432611824951166"; "I"="0.6843865559598825"; "J"="0.2076816565155326"; "K"="1.928934095147895"; "N"="1.157616303759099"; "O"="3.110295588972633" }
    12 = @{ "C"="0.03033499269979245"; "D"="0.2295887022403349"; "E"="0.1918907997309987"; "F"="1.351288292521758"; "G"="0.002431878805951198"; "I"="0.6857215531978298"; "J"="0.2092271947121844"; "K"="1.970908245864734"; "N"="1.1581675518362"; "O"="3.112633183215507" }
    13 = @{ "C"="0.03020231337117707"; "D"="0.2292052282142834"; "E"="0.1915792740015334"; "F"="1.350762460511049"; "G"="0.002432036042829234"; "I"="0.685430897866695"; "J"="0.2088935961075009"; "K"="1.961868731984055"; "N"="1.158045919599047"; "O"="3.112115627007256" }
    14 = @{ "C"="0.0297693927239493"; "D"="0.2279564837349426"; "E"="0.1905654516828221"; "F"="1.349068977446493"; "G"="0.002432551234207639"; "I"="0.6844949723868652"; "J"="0.2078084757465319"; "K"="1.932387497531067"; "N"="1.157660343142553"; "O"="3.110481545654693" }
    15 = @{ "C"="0.02950404419058827"; "D"="0.2271930133006208"; "E"="0.1899460942155002"; "F"="1.348048087046053"; "G"="0.002432868654906477"; "I"="0.6839308817345895"; "J"="0.2071459727549438"; "K"="1.914328342914985"; "N"="1.157432694324825"; "O"="3.109521935478114" }
    16 = @{ "C"="0.02798081671609509"; "D"="0.2228397381242218"; "E"="0.1864219388941422"; "F"="1.342449247309034"; "G"="0.002434716085385808"; "I"="0.6808390435259071"; "J"="0.2033826034944752"; "K"="1.8108242322254"; "N"="1.156259009362458"; "O"="3.104655368911011" }
    17 = @{ "C"="0.02704404354527412"; "D"="0.2201887379278986"; "E"="0.1842824379226755"; "F"="1.339238762878466"; "G"="0.002435874839185548"; "I"="0.67906774017667"; "J"="0.2011034842338972"; "K"="1.747316768417022"; "N"="1.155655962572453"; "O"="3.1022336615261" }
    18 = @{ "C"="0.02650438983305747"; "D"="0.2186711858181809"; "E"="0.1830601166772823"; "F"="1.337474510964412"; "G"="0.002436550685364026"; "I"="0.6780949750061183"; "J"="0.1998034638048694"; "K"="1.710785611487438"; "N"="1.155352166920878"; "O"="3.101047674071793" }
    19 = @{ "C"="0.02632152761800199"; "D"="0.2181586141399805"; "E"="0.1826476799761068"; "F"="1.33689130012543"; "G"="0.002436781125248055"; "I"="0.67777351574248"; "J"="0.1993651655494091"; "K"="1.698416283146514"; "N"="1.155256709837616"; "O"="3.100681631249842" }
    20 = @{ "C"="0.02714385242177286"; "D"="0.2204701935646085"; "E"="0.1845093365061246"; "F"="1.339572001668415"; "G"="0.00243575051941409"; "I"="0.6792515320185544"; "J"="0.2013449754114305"; "K"="1.754077608535226"; "N"="1.155715702981865"; "O"="3.102470035424346" }
    21 = @{ "C"="0.02989658792407113"; "D"="0.2283229741760948"; "E"="0.1908628945832405"; "F"="1.349562972838186"; "G"="0.002432399524425422"; "I"="0.6847679607326427"; "J"="0.208126750785226"; "K"="1.941047072932179"; "N"="1.15777181949008"; "O"="3.110952904107393" }
    22 = @{ "C"="0.03168797393738032"; "D"="0.2335189362827208"; "E"="0.1950886543903678"; "F"="1.356827249865844"; "G"="0.002430292360705425"; "I"="0.6887844682137754"; "J"="0.2126559028557438"; "K"="2.063197476202561"; "N"="1.15949749783897"; "O"="3.118345242007734" }
    23 = @{ "C"="0.03073258788272426"; "D"="0.2307399570411803"; "E"="0.1928265876305986"; "F"="1.352882883494729"; "G"="0.002431409429855869"; "I"="0.6866030991629302"; "J"="0.2102297414330963"; "K"="1.998008353116063"; "N"="1.158541602432265"; "O"="3.114230415204446" }
    24 = @{ "C"="0.02709873222397619"; "D"="0.220342927083351"; "E"="0.1844067316505118"; "F"="1.339421090562013"; "G"="0.002435806694204464"; "I"="0.6791682977573714"; "J"="0.2012357652452579"; "K"="1.751021094441683"; "N"="1.15568856067884"; "O"="3.102362528274114" }
    25 = @{ "C"="0.02314673824931646"; "D"="0.2094049488195964"; "E"="0.1756405719273104"; "F"="1.328047192198682"; "G"="0.002440908057353197"; "I"="0.672908142943399"; "J"="0.1919497417291183"; "K"="1.484485619107772"; "N"="1.154258356547075"; "O"="3.098063831071158" }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = [double]$data[$row][$col]
    }
}